$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trial")

# Update the source parameters - sine wave amplitude (A) and offset (y0)
$ws.Range("C1").Value = 80
$ws.Range("C2").Value = -1143

# Make the "Trial" sheet the active/selected sheet (tabSelected moves here)
$ws.Activate()
